$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = $ws.Range("A13").Value2
$ws.Range("B14").Value = "11. Container With Most Water"
$ws.Range("C14").Value = 'This solution has 2 pointer & has sort of greedy approach. Start with "up" at last position & "down" pointer at indx 0. Calculate the volume by (up-down)*min(height[up], height[down]) and update the ans if its smaller. Now if "up" is at lets say height = 7 , while down is at height = 1 why would u update the ptr "up" to go to a unknown height, insead update the ptr that points to a smaller height in this example "down" ptr. If(height of down < height of up) update down,, else update up,,, edge case if both heights are = it doesnt matter which u update'

$ws.Range("A13:C13").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(14).RowHeight = 72

$ws.Range("C14").Select()
